# adicionado politica de preco
# Insert two new columns ("modelo" and "politica") before the existing
# "full" column, shifting full/tipo/link one column right (C,D,E -> E,F,G),
# then fill in the new columns and refresh the moved link's tracking id.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift C:D (formerly full/tipo) right by two columns, opening up C:D.
$ws.Range("C:D").Insert()

# New header cells.
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"

# New data cells for row 2.
$ws.Range("C2").Value = "Sem Modelo"
$ws.Range("D2").Value = ""

# Data that moved from D2/E2 before the insert now lives in F2/G2 and needs
# its new content (tipo lowercased, link tracking_id refreshed).
$ws.Range("F2").Value = "classico"
$ws.Range("G2").Value = "https://produto.mercadolivre.com.br/MLB-4234174824-processador-jfa-digital-j4-redline-profissional-equalizador-_JM#position%3D2%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D857a0209-f2e8-4e6e-a60e-cc0916c156e6"
